# "test 7 and test 8 update"
# Adds a new "Platform" column (K), shifting the former "Outcome Summary"
# column to L, fixes up a couple of existing cells, and appends two new
# experiment rows (007 and 008) to the "Experiment tracking" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiment tracking")

# --- Move the existing "Outcome Summary" column (K) to the new column L ---
# (copy, including formatting, before K gets overwritten with "Platform" data)
$ws.Range("K1:K7").Copy($ws.Range("L1"))

# --- Populate the new "Platform" column (K) ---
$ws.Range("K1").Value = "Platform"
$ws.Range("K2").Value = "Google Colab"
$ws.Range("K3").Value = "Google Colab"
$ws.Range("K4").Value = "Google Colab"
$ws.Range("K5").Value = "Google Colab"
$ws.Range("K6").Value = "Google Colab"
$ws.Range("K7").Value = "Google Colab"

# --- Fix up H5 (test 4 note) and D7:F7 (test 6 had no pocket constraints) ---
$ws.Range("H5").Value = "added eight ligands"
$ws.Range("D7").Value = "empty"
$ws.Range("E7").Value = "empty"
$ws.Range("F7").Value = "empty"

# --- New row 8: Test 007 ---
$ws.Range("A8").Value = "007"
$ws.Range("B8").Value = "007_fgfr2_ECD_ICD_constraints.yaml"
$ws.Range("C8").Value = "FGFR2_v1"
$ws.Range("D8").Value = "empty"
$ws.Range("E8").Value = "empty"
$ws.Range("F8").Value = "empty"
$ws.Range("G8").Value = "empty"
$ws.Range("H8").Value = "Since no minimum distance is defined in the constraint setup, I set the maximum distance to a very large value (100)."
$ws.Range("I8").Value = "007_fgfr2_ECD_ICD_constraints.cif"
$ws.Range("J8").Value = "N/A"
$ws.Range("K8").Value = "Google Colab"
$ws.Range("L8").Value = "Constraints tend to enforce contact between the two domains rather than separate them."

# --- New row 9: Test 008 ---
$ws.Range("A9").Value = "008"
$ws.Range("B9").Value = "008_fgfr2_ligand_x7.yaml"
$ws.Range("C9").Value = "FGFR2_v1"
$ws.Range("D9").Value = "ligand_v1"
$ws.Range("E9").Value = "A:378, A:398"
$ws.Range("F9").Value = "6.0 Å"
$ws.Range("G9").Value = "empty"
$ws.Range("H9").Value = "Since I could not reproduce Test 4 due to GPU limitations on Google Colab, I tested the maximum number of ligands I can add without using Colab Pro."
$ws.Range("I9").Value = "008_fgfr2_ligand_x7_model.cif"
$ws.Range("J9").Value = "N/A"
$ws.Range("K9").Value = "Google Colab"
$ws.Range("L9").Value = "The maximum number of ligands we can include without a Google Colab Pro subscription is seven."

# --- Column widths (character units; this runtime stores width = input + 5/7) ---
$ws.Columns.Item(8).ColumnWidth = 18.648995535714285   # H -> 19.36328125
$ws.Columns.Item(11).ColumnWidth = 14.375558035714286  # K -> 15.08984375
$ws.Columns.Item(12).ColumnWidth = 17.285714285714285  # L -> 18

# --- Update the active cell selection to match the saved view ---
$ws.Activate()
$ws.Range("D14").Select()
